$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log entries to append after the existing data (rows 1-18)
$newRows = @(
    @("2024-05-13", "11:35:46", "No pone tornillo", "-", "-", "-", "-"),
    @("2024-05-13", "11:35:59", "-", "-", "Detección de sealling mal puesto", "-", "-"),
    @("2024-05-13", "11:54:10", "-", "Cámara no detecta foam derecho", "-", "-", "-")
)

$startRow = 19
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Length; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($c -eq 1) {
            # Force text storage so the date-looking value ("2024-05-13")
            # is not auto-converted into an Excel date serial number.
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowData[$c - 1]
    }
}
